# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value, shared by both sheets that
# contain this dataset.
$updates = @{
    2  = 831
    5  = 50
    6  = 12385
    7  = 54
    9  = 496
    10 = 445
    11 = 1136
    13 = 13621
    14 = 13836
    22 = 241
    23 = 4926
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
